$d = $word.ActiveDocument

$replacements = @(
    @("35×62=", "65×38="),
    @("74×65=", "64×28="),
    @("87×38=", "60×92="),
    @("18×13=", "75×52="),
    @("86×74=", "13×21="),
    @("69×42=", "93×55="),
    @("51×71=", "40×56="),
    @("49×64=", "55×87="),
    @("79×65=", "39×87="),
    @("31×73=", "39×66="),
    @("73×24=", "84×14="),
    @("71×91=", "24×94="),
    @("75×20=", "44×50="),
    @("48×78=", "42×86="),
    @("52×22=", "34×70="),
    @("12×35=", "52×97="),
    @("82×96=", "38×95="),
    @("45×24=", "39×78="),
    @("20×19=", "33×52="),
    @("14×32=", "71×13="),
    @("61×26=", "38×71="),
    @("99×40=", "41×54="),
    @("40×34=", "74×55="),
    @("83×57=", "75×92="),
    @("97×17=", "16×33=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
